# Daily attendance processing - 2025-11-13 06:31:38
#
# Normalizes the "Recorded By" (column G) entries on the active sheet:
# certain exact combinations of recorder names/emails had their
# comma-separated members reordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-text replacements to apply to column G ("Recorded By") values.
$replacements = @{
    "dnasr281@gmail.com, System"              = "System, dnasr281@gmail.com"
    "System, admin@admin.com"                 = "admin@admin.com, System"
    "dnasr281@gmail.com, admin@admin.com"     = "admin@admin.com, dnasr281@gmail.com"
    "backup@backdoor.com, system, System"     = "backup@backdoor.com, System, system"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2
    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value2 = $replacements[$val]
    }
}
